$wb = $excel.ActiveWorkbook

# hunk 0: ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 372.13794
$ws.Range("I92").Value = 349.7143
$ws.Range("K92").Value = 349.7143
$ws.Range("M92").Value = 898.2857

# hunk 1: ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 7633
$ws.Range("I98").Value = 5266.3335
$ws.Range("K98").Value = 5266.3335
$ws.Range("M98").Value = -3768.3335

# hunk 2: ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3566.7334
$ws.Range("I100").Value = 3354.2222
$ws.Range("J100").Value = 3885.5
$ws.Range("K100").Value = 3354.2222
$ws.Range("L100").Value = 3885.5
$ws.Range("M100").Value = -2813.2222
$ws.Range("N100").Value = -4967.5

# hunk 3: ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 800
$ws.Range("I103").Value = 1000
$ws.Range("J103").Value = 700
$ws.Range("K103").Value = 3000
$ws.Range("L103").Value = 2100
$ws.Range("M103").Value = -2414
$ws.Range("N103").Value = -3272

# hunk 4: ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2961.3809
$ws.Range("I116").Value = 2265.5667
$ws.Range("J116").Value = 4700.9165
$ws.Range("K116").Value = 2265.5667
$ws.Range("L116").Value = 4700.9165
$ws.Range("M116").Value = 1176.4333
$ws.Range("N116").Value = -11584.9165

# hunk 5: ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 7633
$ws.Range("I122").Value = 5266.3335
$ws.Range("K122").Value = 15799.0005
$ws.Range("M122").Value = -13349.0005

# hunk 6: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4485.543
$ws.Range("I32").Value = 3054.8572
$ws.Range("K32").Value = 3054.8572
$ws.Range("M32").Value = -2767.8572

# hunk 7: ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1668.1945
$ws.Range("I45").Value = 1088
$ws.Range("K45").Value = 1088
$ws.Range("M45").Value = -711

# hunk 8: ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 866.6
$ws.Range("I110").Value = 521.3570999999999
$ws.Range("K110").Value = 521.3570999999999
$ws.Range("M110").Value = 1523.6429

# hunk 9: BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N94").ClearContents()
$ws.Range("H94").Value = 539.8570999999999
$ws.Range("I94").Value = 539.8570999999999
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 539.8570999999999
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -88.85709999999995

# hunk 10: CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2566820.5
$ws.Range("I31").Value = 3126744
$ws.Range("K31").Value = 3126744
$ws.Range("M31").Value = -3126449

# hunk 11: CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2566820.5
$ws.Range("I34").Value = 3126744
$ws.Range("K34").Value = 3126744
$ws.Range("M34").Value = -3126542

# hunk 12: CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 20836738
$ws.Range("I58").Value = 2529.4285
$ws.Range("J58").Value = 50004630
$ws.Range("K58").Value = 2529.4285
$ws.Range("L58").Value = 50004630
$ws.Range("M58").Value = -2326.4285
$ws.Range("N58").Value = -50005036

# hunk 13: CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2537.625
$ws.Range("I99").Value = 1260.2
$ws.Range("K99").Value = 1260.2
$ws.Range("M99").Value = 237.8

# hunk 14: CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2546.9092
$ws.Range("J107").Value = 3391.9092
$ws.Range("L107").Value = 3391.9092
$ws.Range("N107").Value = -7231.9092

# hunk 15: CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2537.625
$ws.Range("I126").Value = 1260.2
$ws.Range("K126").Value = 3780.6
$ws.Range("M126").Value = -1310.6

# hunk 16: CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 20836738
$ws.Range("I136").Value = 2529.4285
$ws.Range("J136").Value = 50004630
$ws.Range("K136").Value = 7588.2855
$ws.Range("L136").Value = 150013890
$ws.Range("M136").Value = -5038.2855
$ws.Range("N136").Value = -150018990

# hunk 17: CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 164.45
$ws.Range("I12").Value = 13.666667
$ws.Range("J12").Value = 229.07143
$ws.Range("K12").Value = 41.000001
$ws.Range("L12").Value = 687.21429
$ws.Range("M12").Value = 131.999999
$ws.Range("N12").Value = -1033.21429

# hunk 18: CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M98").ClearContents()
$ws.Range("H98").Value = 304
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 304
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 912
$ws.Range("N98").Value = -3908

# hunk 19: CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1875.6
$ws.Range("J107").Value = 2332.6667
$ws.Range("L107").Value = 6998.000100000001
$ws.Range("N107").Value = -10838.0001

# hunk 20: CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3149.4167
$ws.Range("I132").Value = 2899.125
$ws.Range("K132").Value = 26092.125
$ws.Range("M132").Value = -23562.125

# hunk 21: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6875
$ws.Range("I122").Value = 7900
$ws.Range("J122").Value = 6142.857
$ws.Range("K122").Value = 23700
$ws.Range("L122").Value = 18428.571
$ws.Range("M122").Value = -21250
$ws.Range("N122").Value = -23328.571

# hunk 22: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2732.175
$ws.Range("I132").Value = 2084.7778
$ws.Range("J132").Value = 4076.7693
$ws.Range("K132").Value = 6254.3334
$ws.Range("L132").Value = 12230.3079
$ws.Range("M132").Value = -3724.3334
$ws.Range("N132").Value = -17290.3079

# hunk 23: GSM row 137
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 29830
$ws.Range("J137").Value = 29830
$ws.Range("L137").Value = 29830
$ws.Range("N137").Value = -40030

# hunk 24: LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2758.4
$ws.Range("I7").Value = 1919.8
$ws.Range("K7").Value = 1919.8
$ws.Range("M7").Value = -1807.8

# hunk 25: LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2754.0952
$ws.Range("I93").Value = 2090.1177
$ws.Range("J93").Value = 5576
$ws.Range("K93").Value = 2090.1177
$ws.Range("L93").Value = 5576
$ws.Range("M93").Value = -842.1176999999998
$ws.Range("N93").Value = -8072

# hunk 26: LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2758.4
$ws.Range("I126").Value = 1919.8
$ws.Range("K126").Value = 5759.4
$ws.Range("M126").Value = -3289.4

# hunk 27: LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 29745
$ws.Range("J133").Value = 29745
$ws.Range("L133").Value = 29745
$ws.Range("N133").Value = -34805

# hunk 28: LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 29583.334
$ws.Range("J140").Value = 29583.334
$ws.Range("L140").Value = 29583.334
$ws.Range("N140").Value = -39943.334

# hunk 29: WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1680
$ws.Range("I96").Value = 1750
$ws.Range("J96").Value = 1633.3334
$ws.Range("K96").Value = 1750
$ws.Range("L96").Value = 1633.3334
$ws.Range("M96").Value = -377
$ws.Range("N96").Value = -4379.3334

# hunk 30: WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 627298.1
$ws.Range("I122").Value = 668451.3
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 2005353.9
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -2002903.9
$ws.Range("N122").Value = -34900
